# Commit: "Updating ID pattern to accept SenNet ID"
#
# Changes applied:
#  1. Organ sheet, A1 comment: update wording so the ID pattern description
#     covers both HuBMAP and SenNet identifiers.
#  2. weight_unit lookup sheet: add a third allowed unit, "mg", with its
#     ontology term URL.
#  3. Organ sheet data validation on the weight_unit column (M): widen the
#     list range so it includes the new third row.
#  4. .metadata sheet: bump the pav:createdOn timestamp to reflect the
#     new revision.

$wb = $excel.ActiveWorkbook

# --- 1. Update the organ_id comment on the Organ sheet -----------------
$organSheet = $wb.Sheets.Item("Organ")
$newComment = "(Required) Unique HuBMAP or SenNet identifier for the organ. Example:`nHBM811.ORKO.128 or SNT914.IKOK.489"
$null = $organSheet.Range("A1").Comment.Text($newComment)

# --- 2. Add the "mg" unit option to the weight_unit sheet ---------------
$weightUnitSheet = $wb.Sheets.Item("weight_unit")
$weightUnitSheet.Range("A3").Value = "mg"
$weightUnitSheet.Range("B3").Value = "http://purl.obolibrary.org/obo/UO_0000022"

# --- 3. Extend the weight_unit data validation list on column M --------
$organSheet.Range("M2:M1001").Validation.Formula1 = "'weight_unit'!`$A`$1:`$A`$3"

# --- 4. Bump the pav:createdOn timestamp in .metadata -------------------
$metadataSheet = $wb.Sheets.Item(".metadata")
$metadataSheet.Range("C2").Value = "2023-10-27T18:07:03-07:00"
